# previsao_retorno.xlsx - "atualizei dados da bibi e add"
#
# The underlying report was regenerated: every "INATIVO - X meses sem
# comprar" value in the "situacao" column (J) advances by 0.1 month, and
# three "ATIVO" customers (rows 38, 68, 115) got a fresh purchase dated a
# bit later, which shifts their total_compras_historico (E),
# ultima_compra (H) and proxima_compra (I) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) situacao (column J): bump every "INATIVO - X meses sem comprar" by 0.1 ---
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $val = $cell.Value2
    if ($val -is [string] -and $val -match '^INATIVO - ([\d.]+) meses sem comprar$') {
        $num = [double]$matches[1] + 0.1
        $cell.Value = "INATIVO - {0:N1} meses sem comprar" -f $num
    }
}

# --- 2) rows with a new purchase logged: total_compras_historico / ultima_compra / proxima_compra ---
$purchaseUpdates = @(
    @{Row=38;  Total=23;    Ultima=45850.5727662037;  Proxima=45912.5727662037},
    @{Row=68;  Total=33;    Ultima=45850.78201388889; Proxima=45881.78201388889},
    @{Row=115; Total=16551; Ultima=45849.74519675926; Proxima=45850.74519675926}
)
foreach ($u in $purchaseUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Total
    $ws.Cells.Item($u.Row, 8).Value = $u.Ultima
    $ws.Cells.Item($u.Row, 9).Value = $u.Proxima
}
